# Identifying a student by email instead of name.
# The "Ønskede samarbeidspartnere?" (desired collaborators) column (I) used to
# reference students using the short "stdN" form; align it with the
# "Brukernavn" (username) naming scheme used elsewhere in the sheet by
# renaming every "stdN" token to "studentN" (tokens are ';'-separated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 29
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 9)   # column I
    $val = $cell.Value2
    if ($null -eq $val) { continue }
    $text = [string]$val
    if ($text -eq "") { continue }

    $parts = $text -split ";"
    $changed = $false
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i] -match "^std(\d+)$") {
            $parts[$i] = "student" + $matches[1]
            $changed = $true
        }
    }
    if ($changed) {
        $cell.Value = [string]::Join(";", $parts)
    }
}

# Reflect the final state of the interactive editing session: the user ended
# up with cell I29 selected and the view scrolled back to the top of the
# sheet (no frozen/offset top-left cell).
[void]$ws.Range("I29").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
